$wb = $excel.ActiveWorkbook

# LLL_max_6 / LLL_max_10
$wsA = $wb.Worksheets.Item("LLL_max_6")
$wsB = $wb.Worksheets.Item("LLL_max_10")
$wsA.Range("P2").Value = [double]"0"
$wsB.Range("P2").Value = [double]"0"
$wsA.Range("Q2").Value = [double]"0"
$wsB.Range("Q2").Value = [double]"0"
$wsA.Range("P3").Value = [double]"0"
$wsB.Range("P3").Value = [double]"0"
$wsA.Range("Q3").Value = [double]"0"
$wsB.Range("Q3").Value = [double]"0"
$wsA.Range("P4").Value = [double]"0"
$wsB.Range("P4").Value = [double]"0"
$wsA.Range("Q4").Value = [double]"0"
$wsB.Range("Q4").Value = [double]"0"

# LLL_max_fault_6 / LLL_max_fault_10
$wsA = $wb.Worksheets.Item("LLL_max_fault_6")
$wsB = $wb.Worksheets.Item("LLL_max_fault_10")
$wsA.Range("P2").Value = [double]"-10.05101633965329"
$wsB.Range("P2").Value = [double]"-10.05101633965329"
$wsA.Range("Q2").Value = [double]"-9.529673900180386"
$wsB.Range("Q2").Value = [double]"-9.529673900180386"
$wsA.Range("P3").Value = [double]"-9.086152551827228"
$wsB.Range("P3").Value = [double]"-9.086152551827228"
$wsA.Range("Q3").Value = [double]"-9.529673900180386"
$wsB.Range("Q3").Value = [double]"-9.529673900180386"
$wsA.Range("P4").Value = [double]"-9.086152551827231"
$wsB.Range("P4").Value = [double]"-9.086152551827231"
$wsA.Range("Q4").Value = [double]"-9.529673900180386"
$wsB.Range("Q4").Value = [double]"-9.529673900180386"

# LLL_min_6 / LLL_min_10
$wsA = $wb.Worksheets.Item("LLL_min_6")
$wsB = $wb.Worksheets.Item("LLL_min_10")
$wsA.Range("P2").Value = [double]"0"
$wsB.Range("P2").Value = [double]"0"
$wsA.Range("Q2").Value = [double]"0"
$wsB.Range("Q2").Value = [double]"0"
$wsA.Range("P3").Value = [double]"0"
$wsB.Range("P3").Value = [double]"0"
$wsA.Range("Q3").Value = [double]"0"
$wsB.Range("Q3").Value = [double]"0"
$wsA.Range("P4").Value = [double]"0"
$wsB.Range("P4").Value = [double]"0"
$wsA.Range("Q4").Value = [double]"0"
$wsB.Range("Q4").Value = [double]"0"

# LLL_min_fault_6 / LLL_min_fault_10
$wsA = $wb.Worksheets.Item("LLL_min_fault_6")
$wsB = $wb.Worksheets.Item("LLL_min_fault_10")
$wsA.Range("P2").Value = [double]"-10.4045966787725"
$wsB.Range("P2").Value = [double]"-10.4045966787725"
$wsA.Range("Q2").Value = [double]"-7.634345070499261"
$wsB.Range("Q2").Value = [double]"-7.634345070499261"
$wsA.Range("P3").Value = [double]"-5.397652098192608"
$wsB.Range("P3").Value = [double]"-5.397652098192608"
$wsA.Range("Q3").Value = [double]"-7.634345070499261"
$wsB.Range("Q3").Value = [double]"-7.634345070499261"
$wsA.Range("P4").Value = [double]"-5.397652098192621"
$wsB.Range("P4").Value = [double]"-5.397652098192621"
$wsA.Range("Q4").Value = [double]"-7.634345070499261"
$wsB.Range("Q4").Value = [double]"-7.634345070499261"

# LL_max_6 / LL_max_10
$wsA = $wb.Worksheets.Item("LL_max_6")
$wsB = $wb.Worksheets.Item("LL_max_10")
$wsA.Range("AL2").Value = [double]"7.100341988099845E-13"
$wsB.Range("AL2").Value = [double]"7.100341988099845E-13"
$wsA.Range("AM2").Value = [double]"179.9999999999793"
$wsB.Range("AM2").Value = [double]"179.9999999999793"
$wsA.Range("AN2").Value = [double]"-179.9999999999858"
$wsB.Range("AN2").Value = [double]"-179.9999999999858"
$wsA.Range("AO2").Value = [double]"5.400307135557078E-13"
$wsB.Range("AO2").Value = [double]"5.400307135557078E-13"
$wsA.Range("AP2").Value = [double]"179.9999999999542"
$wsB.Range("AP2").Value = [double]"179.9999999999542"
$wsA.Range("AQ2").Value = [double]"-179.9999999999592"
$wsB.Range("AQ2").Value = [double]"-179.9999999999592"
$wsA.Range("AL3").Value = [double]"4.077258211264359E-13"
$wsB.Range("AL3").Value = [double]"4.077258211264359E-13"
$wsA.Range("AM3").Value = [double]"179.9999999999479"
$wsB.Range("AM3").Value = [double]"179.9999999999479"
$wsA.Range("AN3").Value = [double]"-179.9999999999518"
$wsB.Range("AN3").Value = [double]"-179.9999999999518"
$wsA.Range("AO3").Value = [double]"5.400307135557078E-13"
$wsB.Range("AO3").Value = [double]"5.400307135557078E-13"
$wsA.Range("AP3").Value = [double]"179.9999999999542"
$wsB.Range("AP3").Value = [double]"179.9999999999542"
$wsA.Range("AQ3").Value = [double]"-179.9999999999592"
$wsB.Range("AQ3").Value = [double]"-179.9999999999592"
$wsA.Range("AL4").Value = [double]"4.235678099021902E-13"
$wsB.Range("AL4").Value = [double]"4.235678099021902E-13"
$wsA.Range("AM4").Value = [double]"179.9999999999479"
$wsB.Range("AM4").Value = [double]"179.9999999999479"
$wsA.Range("AN4").Value = [double]"-179.9999999999519"
$wsB.Range("AN4").Value = [double]"-179.9999999999519"
$wsA.Range("AO4").Value = [double]"5.400307135557078E-13"
$wsB.Range("AO4").Value = [double]"5.400307135557078E-13"
$wsA.Range("AP4").Value = [double]"179.9999999999542"
$wsB.Range("AP4").Value = [double]"179.9999999999542"
$wsA.Range("AQ4").Value = [double]"-179.9999999999592"
$wsB.Range("AQ4").Value = [double]"-179.9999999999592"

# LL_max_fault_6 / LL_max_fault_10
$wsA = $wb.Worksheets.Item("LL_max_fault_6")
$wsB = $wb.Worksheets.Item("LL_max_fault_10")
$wsA.Range("AL2").Value = [double]"2.594692874650617E-13"
$wsB.Range("AL2").Value = [double]"2.594692874650617E-13"
$wsA.Range("AM2").Value = [double]"-140.8543229646633"
$wsB.Range("AM2").Value = [double]"-140.8543229646633"
$wsA.Range("AN2").Value = [double]"122.4159013159207"
$wsB.Range("AN2").Value = [double]"122.4159013159207"
$wsA.Range("AO2").Value = [double]"2.191965202687454E-13"
$wsB.Range("AO2").Value = [double]"2.191965202687454E-13"
$wsA.Range("AP2").Value = [double]"-142.8672784204022"
$wsB.Range("AP2").Value = [double]"-142.8672784204022"
$wsA.Range("AQ2").Value = [double]"128.2389811279577"
$wsB.Range("AQ2").Value = [double]"128.2389811279577"
$wsA.Range("AL3").Value = [double]"1.866433575455782E-13"
$wsB.Range("AL3").Value = [double]"1.866433575455782E-13"
$wsA.Range("AM3").Value = [double]"-144.8695745192862"
$wsB.Range("AM3").Value = [double]"-144.8695745192862"
$wsA.Range("AN3").Value = [double]"133.0893358579405"
$wsB.Range("AN3").Value = [double]"133.0893358579405"
$wsA.Range("AO3").Value = [double]"2.191965202687454E-13"
$wsB.Range("AO3").Value = [double]"2.191965202687454E-13"
$wsA.Range("AP3").Value = [double]"-142.8672784204022"
$wsB.Range("AP3").Value = [double]"-142.8672784204022"
$wsA.Range("AQ3").Value = [double]"128.2389811279577"
$wsB.Range("AQ3").Value = [double]"128.2389811279577"
$wsA.Range("AL4").Value = [double]"1.920934498077301E-13"
$wsB.Range("AL4").Value = [double]"1.920934498077301E-13"
$wsA.Range("AM4").Value = [double]"-144.8695745192862"
$wsB.Range("AM4").Value = [double]"-144.8695745192862"
$wsA.Range("AN4").Value = [double]"133.0893358579404"
$wsB.Range("AN4").Value = [double]"133.0893358579404"
$wsA.Range("AO4").Value = [double]"2.191965202687454E-13"
$wsB.Range("AO4").Value = [double]"2.191965202687454E-13"
$wsA.Range("AP4").Value = [double]"-142.8672784204022"
$wsB.Range("AP4").Value = [double]"-142.8672784204022"
$wsA.Range("AQ4").Value = [double]"128.2389811279577"
$wsB.Range("AQ4").Value = [double]"128.2389811279577"

# LL_min_6 / LL_min_10
$wsA = $wb.Worksheets.Item("LL_min_6")
$wsB = $wb.Worksheets.Item("LL_min_10")
$wsA.Range("AL2").Value = [double]"7.2703997191632E-13"
$wsB.Range("AL2").Value = [double]"7.2703997191632E-13"
$wsA.Range("AM2").Value = [double]"-179.9999999999019"
$wsB.Range("AM2").Value = [double]"-179.9999999999019"
$wsA.Range("AN2").Value = [double]"179.9999999998949"
$wsB.Range("AN2").Value = [double]"179.9999999998949"
$wsA.Range("AO2").Value = [double]"4.836847457647291E-13"
$wsB.Range("AO2").Value = [double]"4.836847457647291E-13"
$wsA.Range("AP2").Value = [double]"-179.9999999994768"
$wsB.Range("AP2").Value = [double]"-179.9999999994768"
$wsA.Range("AQ2").Value = [double]"179.999999999472"
$wsB.Range("AQ2").Value = [double]"179.999999999472"
$wsA.Range("AL3").Value = [double]"3.304768388060126E-13"
$wsB.Range("AL3").Value = [double]"3.304768388060126E-13"
$wsA.Range("AM3").Value = [double]"-179.9999999994135"
$wsB.Range("AM3").Value = [double]"-179.9999999994135"
$wsA.Range("AN3").Value = [double]"179.99999999941"
$wsB.Range("AN3").Value = [double]"179.99999999941"
$wsA.Range("AO3").Value = [double]"4.836847457647291E-13"
$wsB.Range("AO3").Value = [double]"4.836847457647291E-13"
$wsA.Range("AP3").Value = [double]"-179.9999999994768"
$wsB.Range("AP3").Value = [double]"-179.9999999994768"
$wsA.Range("AQ3").Value = [double]"179.999999999472"
$wsB.Range("AQ3").Value = [double]"179.999999999472"
$wsA.Range("AL4").Value = [double]"3.307476491316185E-13"
$wsB.Range("AL4").Value = [double]"3.307476491316185E-13"
$wsA.Range("AM4").Value = [double]"-179.9999999994134"
$wsB.Range("AM4").Value = [double]"-179.9999999994134"
$wsA.Range("AN4").Value = [double]"179.99999999941"
$wsB.Range("AN4").Value = [double]"179.99999999941"
$wsA.Range("AO4").Value = [double]"4.836847457647291E-13"
$wsB.Range("AO4").Value = [double]"4.836847457647291E-13"
$wsA.Range("AP4").Value = [double]"-179.9999999994768"
$wsB.Range("AP4").Value = [double]"-179.9999999994768"
$wsA.Range("AQ4").Value = [double]"179.999999999472"
$wsB.Range("AQ4").Value = [double]"179.999999999472"

# LL_min_fault_6 / LL_min_fault_10
$wsA = $wb.Worksheets.Item("LL_min_fault_6")
$wsB = $wb.Worksheets.Item("LL_min_fault_10")
$wsA.Range("AL2").Value = [double]"2.90237442368646E-13"
$wsB.Range("AL2").Value = [double]"2.90237442368646E-13"
$wsA.Range("AM2").Value = [double]"-141.4823925534191"
$wsB.Range("AM2").Value = [double]"-141.4823925534191"
$wsA.Range("AN2").Value = [double]"122.8363393546082"
$wsB.Range("AN2").Value = [double]"122.8363393546082"
$wsA.Range("AO2").Value = [double]"2.15386192807484E-13"
$wsB.Range("AO2").Value = [double]"2.15386192807484E-13"
$wsA.Range("AP2").Value = [double]"-142.5584264979142"
$wsB.Range("AP2").Value = [double]"-142.5584264979142"
$wsA.Range("AQ2").Value = [double]"131.6989175372502"
$wsB.Range("AQ2").Value = [double]"131.6989175372502"
$wsA.Range("AL3").Value = [double]"1.628916501490903E-13"
$wsB.Range("AL3").Value = [double]"1.628916501490903E-13"
$wsA.Range("AM3").Value = [double]"-144.4405405825464"
$wsB.Range("AM3").Value = [double]"-144.4405405825464"
$wsA.Range("AN3").Value = [double]"138.2820504968733"
$wsB.Range("AN3").Value = [double]"138.2820504968733"
$wsA.Range("AO3").Value = [double]"2.15386192807484E-13"
$wsB.Range("AO3").Value = [double]"2.15386192807484E-13"
$wsA.Range("AP3").Value = [double]"-142.5584264979142"
$wsB.Range("AP3").Value = [double]"-142.5584264979142"
$wsA.Range("AQ3").Value = [double]"131.6989175372502"
$wsB.Range("AQ3").Value = [double]"131.6989175372502"
$wsA.Range("AL4").Value = [double]"1.676833193653677E-13"
$wsB.Range("AL4").Value = [double]"1.676833193653677E-13"
$wsA.Range("AM4").Value = [double]"-144.4405405825464"
$wsB.Range("AM4").Value = [double]"-144.4405405825464"
$wsA.Range("AN4").Value = [double]"138.2820504968733"
$wsB.Range("AN4").Value = [double]"138.2820504968733"
$wsA.Range("AO4").Value = [double]"2.15386192807484E-13"
$wsB.Range("AO4").Value = [double]"2.15386192807484E-13"
$wsA.Range("AP4").Value = [double]"-142.5584264979142"
$wsB.Range("AP4").Value = [double]"-142.5584264979142"
$wsA.Range("AQ4").Value = [double]"131.6989175372502"
$wsB.Range("AQ4").Value = [double]"131.6989175372502"

# LLG_max_6 / LLG_max_10
$wsA = $wb.Worksheets.Item("LLG_max_6")
$wsB = $wb.Worksheets.Item("LLG_max_10")
$wsA.Range("AL2").Value = [double]"-0.02683261092834746"
$wsB.Range("AL2").Value = [double]"-0.02683261092834746"
$wsA.Range("AM2").Value = [double]"0"
$wsB.Range("AM2").Value = [double]"0"
$wsA.Range("AN2").Value = [double]"0"
$wsB.Range("AN2").Value = [double]"0"
$wsA.Range("AO2").Value = [double]"-10.99131251259726"
$wsB.Range("AO2").Value = [double]"-10.99131251259726"
$wsA.Range("AP2").Value = [double]"0"
$wsB.Range("AP2").Value = [double]"0"
$wsA.Range("AQ2").Value = [double]"0"
$wsB.Range("AQ2").Value = [double]"0"
$wsA.Range("AL3").Value = [double]"-7.597588372061134"
$wsB.Range("AL3").Value = [double]"-7.597588372061134"
$wsA.Range("AM3").Value = [double]"0"
$wsB.Range("AM3").Value = [double]"0"
$wsA.Range("AN3").Value = [double]"0"
$wsB.Range("AN3").Value = [double]"0"
$wsA.Range("AO3").Value = [double]"-10.99131251259726"
$wsB.Range("AO3").Value = [double]"-10.99131251259726"
$wsA.Range("AP3").Value = [double]"0"
$wsB.Range("AP3").Value = [double]"0"
$wsA.Range("AQ3").Value = [double]"0"
$wsB.Range("AQ3").Value = [double]"0"
$wsA.Range("AL4").Value = [double]"-7.59758837206113"
$wsB.Range("AL4").Value = [double]"-7.59758837206113"
$wsA.Range("AM4").Value = [double]"0"
$wsB.Range("AM4").Value = [double]"0"
$wsA.Range("AN4").Value = [double]"0"
$wsB.Range("AN4").Value = [double]"0"
$wsA.Range("AO4").Value = [double]"-10.99131251259726"
$wsB.Range("AO4").Value = [double]"-10.99131251259726"
$wsA.Range("AP4").Value = [double]"0"
$wsB.Range("AP4").Value = [double]"0"
$wsA.Range("AQ4").Value = [double]"0"
$wsB.Range("AQ4").Value = [double]"0"

# LLG_max_fault_6 / LLG_max_fault_10
$wsA = $wb.Worksheets.Item("LLG_max_fault_6")
$wsB = $wb.Worksheets.Item("LLG_max_fault_10")
$wsA.Range("AL2").Value = [double]"-2.047021704634401"
$wsB.Range("AL2").Value = [double]"-2.047021704634401"
$wsA.Range("AM2").Value = [double]"-131.3437653417182"
$wsB.Range("AM2").Value = [double]"-131.3437653417182"
$wsA.Range("AN2").Value = [double]"113.278029361874"
$wsB.Range("AN2").Value = [double]"113.278029361874"
$wsA.Range("AO2").Value = [double]"-4.308903017906394"
$wsB.Range("AO2").Value = [double]"-4.308903017906394"
$wsA.Range("AP2").Value = [double]"-124.6110161336328"
$wsB.Range("AP2").Value = [double]"-124.6110161336328"
$wsA.Range("AQ2").Value = [double]"109.8444437755527"
$wsB.Range("AQ2").Value = [double]"109.8444437755527"
$wsA.Range("AL3").Value = [double]"-4.775555041063019"
$wsB.Range("AL3").Value = [double]"-4.775555041063019"
$wsA.Range("AM3").Value = [double]"-120.4649322159348"
$wsB.Range("AM3").Value = [double]"-120.4649322159348"
$wsA.Range("AN3").Value = [double]"106.8663009952652"
$wsB.Range("AN3").Value = [double]"106.8663009952652"
$wsA.Range("AO3").Value = [double]"-4.308903017906394"
$wsB.Range("AO3").Value = [double]"-4.308903017906394"
$wsA.Range("AP3").Value = [double]"-124.6110161336328"
$wsB.Range("AP3").Value = [double]"-124.6110161336328"
$wsA.Range("AQ3").Value = [double]"109.8444437755527"
$wsB.Range("AQ3").Value = [double]"109.8444437755527"
$wsA.Range("AL4").Value = [double]"-4.775555041063012"
$wsB.Range("AL4").Value = [double]"-4.775555041063012"
$wsA.Range("AM4").Value = [double]"-120.4649322159348"
$wsB.Range("AM4").Value = [double]"-120.4649322159348"
$wsA.Range("AN4").Value = [double]"106.8663009952652"
$wsB.Range("AN4").Value = [double]"106.8663009952652"
$wsA.Range("AO4").Value = [double]"-4.308903017906394"
$wsB.Range("AO4").Value = [double]"-4.308903017906394"
$wsA.Range("AP4").Value = [double]"-124.6110161336328"
$wsB.Range("AP4").Value = [double]"-124.6110161336328"
$wsA.Range("AQ4").Value = [double]"109.8444437755527"
$wsB.Range("AQ4").Value = [double]"109.8444437755527"

# LLG_min_6 / LLG_min_10
$wsA = $wb.Worksheets.Item("LLG_min_6")
$wsB = $wb.Worksheets.Item("LLG_min_10")
$wsA.Range("AL2").Value = [double]"-0.07043074203077801"
$wsB.Range("AL2").Value = [double]"-0.07043074203077801"
$wsA.Range("AM2").Value = [double]"0"
$wsB.Range("AM2").Value = [double]"0"
$wsA.Range("AN2").Value = [double]"0"
$wsB.Range("AN2").Value = [double]"0"
$wsA.Range("AO2").Value = [double]"-8.312073893397944"
$wsB.Range("AO2").Value = [double]"-8.312073893397944"
$wsA.Range("AP2").Value = [double]"0"
$wsB.Range("AP2").Value = [double]"0"
$wsA.Range("AQ2").Value = [double]"0"
$wsB.Range("AQ2").Value = [double]"0"
$wsA.Range("AL3").Value = [double]"-5.141094698820867"
$wsB.Range("AL3").Value = [double]"-5.141094698820867"
$wsA.Range("AM3").Value = [double]"0"
$wsB.Range("AM3").Value = [double]"0"
$wsA.Range("AN3").Value = [double]"0"
$wsB.Range("AN3").Value = [double]"0"
$wsA.Range("AO3").Value = [double]"-8.312073893397944"
$wsB.Range("AO3").Value = [double]"-8.312073893397944"
$wsA.Range("AP3").Value = [double]"0"
$wsB.Range("AP3").Value = [double]"0"
$wsA.Range("AQ3").Value = [double]"0"
$wsB.Range("AQ3").Value = [double]"0"
$wsA.Range("AL4").Value = [double]"-5.141094698820842"
$wsB.Range("AL4").Value = [double]"-5.141094698820842"
$wsA.Range("AM4").Value = [double]"0"
$wsB.Range("AM4").Value = [double]"0"
$wsA.Range("AN4").Value = [double]"0"
$wsB.Range("AN4").Value = [double]"0"
$wsA.Range("AO4").Value = [double]"-8.312073893397944"
$wsB.Range("AO4").Value = [double]"-8.312073893397944"
$wsA.Range("AP4").Value = [double]"0"
$wsB.Range("AP4").Value = [double]"0"
$wsA.Range("AQ4").Value = [double]"0"
$wsB.Range("AQ4").Value = [double]"0"

# LLG_min_fault_6 / LLG_min_fault_10
$wsA = $wb.Worksheets.Item("LLG_min_fault_6")
$wsB = $wb.Worksheets.Item("LLG_min_fault_10")
$wsA.Range("AL2").Value = [double]"-2.070858898527459"
$wsB.Range("AL2").Value = [double]"-2.070858898527459"
$wsA.Range("AM2").Value = [double]"-131.7373295207326"
$wsB.Range("AM2").Value = [double]"-131.7373295207326"
$wsA.Range("AN2").Value = [double]"112.9875977716212"
$wsB.Range("AN2").Value = [double]"112.9875977716212"
$wsA.Range("AO2").Value = [double]"-5.549016348167304"
$wsB.Range("AO2").Value = [double]"-5.549016348167304"
$wsA.Range("AP2").Value = [double]"-119.0419281637489"
$wsB.Range("AP2").Value = [double]"-119.0419281637489"
$wsA.Range("AQ2").Value = [double]"109.1357001831912"
$wsB.Range("AQ2").Value = [double]"109.1357001831912"
$wsA.Range("AL3").Value = [double]"-5.292660328418042"
$wsB.Range("AL3").Value = [double]"-5.292660328418042"
$wsA.Range("AM3").Value = [double]"-112.3238428218692"
$wsB.Range("AM3").Value = [double]"-112.3238428218692"
$wsA.Range("AN3").Value = [double]"106.1154176972715"
$wsB.Range("AN3").Value = [double]"106.1154176972715"
$wsA.Range("AO3").Value = [double]"-5.549016348167304"
$wsB.Range("AO3").Value = [double]"-5.549016348167304"
$wsA.Range("AP3").Value = [double]"-119.0419281637489"
$wsB.Range("AP3").Value = [double]"-119.0419281637489"
$wsA.Range("AQ3").Value = [double]"109.1357001831912"
$wsB.Range("AQ3").Value = [double]"109.1357001831912"
$wsA.Range("AL4").Value = [double]"-5.292660328418039"
$wsB.Range("AL4").Value = [double]"-5.292660328418039"
$wsA.Range("AM4").Value = [double]"-112.3238428218692"
$wsB.Range("AM4").Value = [double]"-112.3238428218692"
$wsA.Range("AN4").Value = [double]"106.1154176972715"
$wsB.Range("AN4").Value = [double]"106.1154176972715"
$wsA.Range("AO4").Value = [double]"-5.549016348167304"
$wsB.Range("AO4").Value = [double]"-5.549016348167304"
$wsA.Range("AP4").Value = [double]"-119.0419281637489"
$wsB.Range("AP4").Value = [double]"-119.0419281637489"
$wsA.Range("AQ4").Value = [double]"109.1357001831912"
$wsB.Range("AQ4").Value = [double]"109.1357001831912"

# LG_max_6 / LG_max_10
$wsA = $wb.Worksheets.Item("LG_max_6")
$wsB = $wb.Worksheets.Item("LG_max_10")
$wsA.Range("AL2").Value = [double]"0"
$wsB.Range("AL2").Value = [double]"0"
$wsA.Range("AM2").Value = [double]"-105.6656179757978"
$wsB.Range("AM2").Value = [double]"-105.6656179757978"
$wsA.Range("AN2").Value = [double]"105.6598069806291"
$wsB.Range("AN2").Value = [double]"105.6598069806291"
$wsA.Range("AO2").Value = [double]"0"
$wsB.Range("AO2").Value = [double]"0"
$wsA.Range("AP2").Value = [double]"-138.1579831476276"
$wsB.Range("AP2").Value = [double]"-138.1579831476276"
$wsA.Range("AQ2").Value = [double]"117.9292480532439"
$wsB.Range("AQ2").Value = [double]"117.9292480532439"
$wsA.Range("AL3").Value = [double]"0"
$wsB.Range("AL3").Value = [double]"0"
$wsA.Range("AM3").Value = [double]"-145.3006619393671"
$wsB.Range("AM3").Value = [double]"-145.3006619393671"
$wsA.Range("AN3").Value = [double]"124.5503549503629"
$wsB.Range("AN3").Value = [double]"124.5503549503629"
$wsA.Range("AO3").Value = [double]"0"
$wsB.Range("AO3").Value = [double]"0"
$wsA.Range("AP3").Value = [double]"-138.1579831476276"
$wsB.Range("AP3").Value = [double]"-138.1579831476276"
$wsA.Range("AQ3").Value = [double]"117.9292480532439"
$wsB.Range("AQ3").Value = [double]"117.9292480532439"
$wsA.Range("AL4").Value = [double]"0"
$wsB.Range("AL4").Value = [double]"0"
$wsA.Range("AM4").Value = [double]"-145.3006619393671"
$wsB.Range("AM4").Value = [double]"-145.3006619393671"
$wsA.Range("AN4").Value = [double]"124.5503549503629"
$wsB.Range("AN4").Value = [double]"124.5503549503629"
$wsA.Range("AO4").Value = [double]"0"
$wsB.Range("AO4").Value = [double]"0"
$wsA.Range("AP4").Value = [double]"-138.1579831476276"
$wsB.Range("AP4").Value = [double]"-138.1579831476276"
$wsA.Range("AQ4").Value = [double]"117.9292480532439"
$wsB.Range("AQ4").Value = [double]"117.9292480532439"

# LG_max_fault_6 / LG_max_fault_10
$wsA = $wb.Worksheets.Item("LG_max_fault_6")
$wsB = $wb.Worksheets.Item("LG_max_fault_10")
$wsA.Range("AL2").Value = [double]"-8.351765333507215"
$wsB.Range("AL2").Value = [double]"-8.351765333507215"
$wsA.Range("AM2").Value = [double]"-118.5339422527549"
$wsB.Range("AM2").Value = [double]"-118.5339422527549"
$wsA.Range("AN2").Value = [double]"116.8442630436326"
$wsB.Range("AN2").Value = [double]"116.8442630436326"
$wsA.Range("AO2").Value = [double]"-4.157054826293857"
$wsB.Range("AO2").Value = [double]"-4.157054826293857"
$wsA.Range("AP2").Value = [double]"-125.4570209785861"
$wsB.Range("AP2").Value = [double]"-125.4570209785861"
$wsA.Range("AQ2").Value = [double]"120.0875191122162"
$wsB.Range("AQ2").Value = [double]"120.0875191122162"
$wsA.Range("AL3").Value = [double]"-1.174854820720038"
$wsB.Range("AL3").Value = [double]"-1.174854820720038"
$wsA.Range("AM3").Value = [double]"-130.1685625970616"
$wsB.Range("AM3").Value = [double]"-130.1685625970616"
$wsA.Range("AN3").Value = [double]"122.5002041917308"
$wsB.Range("AN3").Value = [double]"122.5002041917308"
$wsA.Range("AO3").Value = [double]"-4.157054826293857"
$wsB.Range("AO3").Value = [double]"-4.157054826293857"
$wsA.Range("AP3").Value = [double]"-125.4570209785861"
$wsB.Range("AP3").Value = [double]"-125.4570209785861"
$wsA.Range("AQ3").Value = [double]"120.0875191122162"
$wsB.Range("AQ3").Value = [double]"120.0875191122162"
$wsA.Range("AL4").Value = [double]"-1.174854820720031"
$wsB.Range("AL4").Value = [double]"-1.174854820720031"
$wsA.Range("AM4").Value = [double]"-130.1685625970616"
$wsB.Range("AM4").Value = [double]"-130.1685625970616"
$wsA.Range("AN4").Value = [double]"122.5002041917308"
$wsB.Range("AN4").Value = [double]"122.5002041917308"
$wsA.Range("AO4").Value = [double]"-4.157054826293857"
$wsB.Range("AO4").Value = [double]"-4.157054826293857"
$wsA.Range("AP4").Value = [double]"-125.4570209785861"
$wsB.Range("AP4").Value = [double]"-125.4570209785861"
$wsA.Range("AQ4").Value = [double]"120.0875191122162"
$wsB.Range("AQ4").Value = [double]"120.0875191122162"

# LG_min_6 / LG_min_10
$wsA = $wb.Worksheets.Item("LG_min_6")
$wsB = $wb.Worksheets.Item("LG_min_10")
$wsA.Range("AL2").Value = [double]"0"
$wsB.Range("AL2").Value = [double]"0"
$wsA.Range("AM2").Value = [double]"-106.0394703642315"
$wsB.Range("AM2").Value = [double]"-106.0394703642315"
$wsA.Range("AN2").Value = [double]"106.023380610698"
$wsB.Range("AN2").Value = [double]"106.023380610698"
$wsA.Range("AO2").Value = [double]"0"
$wsB.Range("AO2").Value = [double]"0"
$wsA.Range("AP2").Value = [double]"-148.3705990935352"
$wsB.Range("AP2").Value = [double]"-148.3705990935352"
$wsA.Range("AQ2").Value = [double]"124.0256852899936"
$wsB.Range("AQ2").Value = [double]"124.0256852899936"
$wsA.Range("AL3").Value = [double]"0"
$wsB.Range("AL3").Value = [double]"0"
$wsA.Range("AM3").Value = [double]"-149.4085172538034"
$wsB.Range("AM3").Value = [double]"-149.4085172538034"
$wsA.Range("AN3").Value = [double]"130.8537654848539"
$wsB.Range("AN3").Value = [double]"130.8537654848539"
$wsA.Range("AO3").Value = [double]"0"
$wsB.Range("AO3").Value = [double]"0"
$wsA.Range("AP3").Value = [double]"-148.3705990935352"
$wsB.Range("AP3").Value = [double]"-148.3705990935352"
$wsA.Range("AQ3").Value = [double]"124.0256852899936"
$wsB.Range("AQ3").Value = [double]"124.0256852899936"
$wsA.Range("AL4").Value = [double]"0"
$wsB.Range("AL4").Value = [double]"0"
$wsA.Range("AM4").Value = [double]"-149.4085172538034"
$wsB.Range("AM4").Value = [double]"-149.4085172538034"
$wsA.Range("AN4").Value = [double]"130.8537654848539"
$wsB.Range("AN4").Value = [double]"130.8537654848539"
$wsA.Range("AO4").Value = [double]"0"
$wsB.Range("AO4").Value = [double]"0"
$wsA.Range("AP4").Value = [double]"-148.3705990935352"
$wsB.Range("AP4").Value = [double]"-148.3705990935352"
$wsA.Range("AQ4").Value = [double]"124.0256852899936"
$wsB.Range("AQ4").Value = [double]"124.0256852899936"

# LG_min_fault_6 / LG_min_fault_10
$wsA = $wb.Worksheets.Item("LG_min_fault_6")
$wsB = $wb.Worksheets.Item("LG_min_fault_10")
$wsA.Range("AL2").Value = [double]"-8.692698158235416"
$wsB.Range("AL2").Value = [double]"-8.692698158235416"
$wsA.Range("AM2").Value = [double]"-118.4884237755319"
$wsB.Range("AM2").Value = [double]"-118.4884237755319"
$wsA.Range("AN2").Value = [double]"116.7866489320772"
$wsB.Range("AN2").Value = [double]"116.7866489320772"
$wsA.Range("AO2").Value = [double]"1.213531850035941"
$wsB.Range("AO2").Value = [double]"1.213531850035941"
$wsA.Range("AP2").Value = [double]"-130.2490184990389"
$wsB.Range("AP2").Value = [double]"-130.2490184990389"
$wsA.Range("AQ2").Value = [double]"121.5733531540493"
$wsB.Range("AQ2").Value = [double]"121.5733531540493"
$wsA.Range("AL3").Value = [double]"7.181611621190969"
$wsB.Range("AL3").Value = [double]"7.181611621190969"
$wsA.Range("AM3").Value = [double]"-136.4394934069981"
$wsB.Range("AM3").Value = [double]"-136.4394934069981"
$wsA.Range("AN3").Value = [double]"125.12124607394"
$wsB.Range("AN3").Value = [double]"125.12124607394"
$wsA.Range("AO3").Value = [double]"1.213531850035941"
$wsB.Range("AO3").Value = [double]"1.213531850035941"
$wsA.Range("AP3").Value = [double]"-130.2490184990389"
$wsB.Range("AP3").Value = [double]"-130.2490184990389"
$wsA.Range("AQ3").Value = [double]"121.5733531540493"
$wsB.Range("AQ3").Value = [double]"121.5733531540493"
$wsA.Range("AL4").Value = [double]"7.181611621190973"
$wsB.Range("AL4").Value = [double]"7.181611621190973"
$wsA.Range("AM4").Value = [double]"-136.4394934069981"
$wsB.Range("AM4").Value = [double]"-136.4394934069981"
$wsA.Range("AN4").Value = [double]"125.1212460739401"
$wsB.Range("AN4").Value = [double]"125.1212460739401"
$wsA.Range("AO4").Value = [double]"1.213531850035941"
$wsB.Range("AO4").Value = [double]"1.213531850035941"
$wsA.Range("AP4").Value = [double]"-130.2490184990389"
$wsB.Range("AP4").Value = [double]"-130.2490184990389"
$wsA.Range("AQ4").Value = [double]"121.5733531540493"
$wsB.Range("AQ4").Value = [double]"121.5733531540493"

Write-Output "done"
